$d = $word.ActiveDocument

# The document contains a YouTube timestamp reference "1:28:33" that needs
# to be updated to "1:43:14" (runs "2"->"43", "8"->":", ":"->"14", and the
# trailing "33" run removed, per the target revision).
$d.Content.Find.Execute("1:28:33", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1:43:14", 2)
